$d = $word.ActiveDocument

# --- 1. Collapse the hyperlink's three runs ("https://github.com/Yu" + "g" + "Kotak/SPM")
#        into a single run, preserving the Hyperlink character style / bold formatting. ---
$hl = $d.Hyperlinks.Item(1)
$hlRange = $d.Range($hl.Range.Start, $hl.Range.End)
$hlRange.Find.Execute("https://github.com/YugKotak/SPM", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "https://github.com/YugKotak/SPM", 2)

# Re-apply the Hyperlink character style to the (now merged) run so it keeps <w:rStyle w:val="Hyperlink"/>
$hl2 = $d.Hyperlinks.Item(1)
$hlRange2 = $d.Range($hl2.Range.Start, $hl2.Range.End)
$hlRange2.Style = "Hyperlink"

# --- 2. "Week 3:" -> "Week 4:" (only the run containing the digit changes) ---
$wkStart = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -ge 7 -and $t.Substring(0, 7) -eq "Week 3:") {
        $wkStart = $p.Range.Start
    }
}
$digitRange = $d.Range($wkStart + 5, $wkStart + 6)
$digitRange.Find.Execute("3", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "4", 2)
